# Swap the data-row values between row 2 and row 3 for columns
# D (Fecha), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado), R (Origen), S (Precio $/Kg).
# All other columns (A,B,C,E-L,Q,T) are identical between the two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "M", "N", "O", "P", "R", "S")

foreach ($col in $cols) {
    $cellRow2 = $ws.Range($col + "2")
    $cellRow3 = $ws.Range($col + "3")

    $val2 = $cellRow2.Value2
    $val3 = $cellRow3.Value2

    $cellRow2.Value2 = $val3
    $cellRow3.Value2 = $val2
}
